$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data for rows 2 and 3 (previously holding the "Quillota" week)
# must be swapped with rows 4 and 5 (previously holding the "O'Higgins" week).
# Only columns D, L, M, N, O, P, R, S differ between the two weeks; the rest
# stay the same for every row, so we only need to swap those columns.

$cols = @("D", "L", "M", "N", "O", "P", "R", "S")

foreach ($col in $cols) {
    $topCell = $ws.Range($col + "2")
    $botCell = $ws.Range($col + "4")
    $tmp = $topCell.Value2
    $topCell.Value = $botCell.Value2
    $botCell.Value = $tmp
}

foreach ($col in $cols) {
    $topCell = $ws.Range($col + "3")
    $botCell = $ws.Range($col + "5")
    $tmp = $topCell.Value2
    $topCell.Value = $botCell.Value2
    $botCell.Value = $tmp
}
